$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update unit-price column (C) values across several rows
$ws.Range("C4").Value = 2800
$ws.Range("C5").Value = 18000
$ws.Range("C6").Value = 1500
$ws.Range("C7").Value = 1500
$ws.Range("C8").Value = 2500
$ws.Range("C24").Value = 2000
$ws.Range("C25").Value = 3500
$ws.Range("C28").Value = 18000
$ws.Range("C33").Value = 6000
$ws.Range("C34").Value = 6000
$ws.Range("C35").Value = 6000
$ws.Range("C36").Value = 6000

# Move the active selection to C28 to match the saved view state
$ws.Range("C28").Select()
